$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1: copy formatting from the existing "sum" header (G1),
# which carries the bold/border/center style, then set its text to "Save".
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New "Save" data column H2:H9, all zeros (era data placeholder).
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
